$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Household"

# Add new data rows
$ws.Range("A2").Value = 101
$ws.Range("B2").Value = "50 000"
$ws.Range("C2").Value = "f"

$ws.Range("A3").Value = 102
$ws.Range("B3").Value = 45000
$ws.Range("C3").Value = "m"

$ws.Range("A4").Value = 103
$ws.Range("B4").Value = 78000
$ws.Range("C4").Value = "f"

# Update selection to match target state
$ws.Range("C5").Select()
